$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher gains a value
$ws.Range("B9").Value = "Alvearie Team"

# Duplicate "Contact" rows become Jurisdiction + (shifted up) Description
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Codes indicating the source of legal document used for identification"

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# Former row 15 (old Immutable/BooleanType[null]) is now redundant - remove it
$ws.Rows.Item(15).Delete()
